$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F5").Value = 9149
$wsExhibit.Range("F6").Value = 9149
$wsExhibit.Range("F7").Value = 535
$wsExhibit.Range("F10").Value = 220
$wsExhibit.Range("F13").Value = 136
$wsExhibit.Range("F16").Value = 11704
$wsExhibit.Range("F17").Value = 11704
$wsExhibit.Range("F23").Value = 391
$wsExhibit.Range("F24").Value = 219
$wsExhibit.Range("F27").Value = 168
$wsExhibit.Range("F28").Value = 145
$wsExhibit.Range("F32").Value = 2091
$wsExhibit.Range("F33").Value = 57
$wsExhibit.Range("F35").Value = 2132
$wsExhibit.Range("F36").Value = 961
$wsExhibit.Range("F37").Value = 4169
$wsExhibit.Range("F39").Value = 3589
$wsExhibit.Range("F40").Value = 329
$wsExhibit.Range("F41").Value = 2606
$wsExhibit.Range("F43").Value = 1295
$wsExhibit.Range("F45").Value = 766
$wsExhibit.Range("F46").Value = 388
$wsExhibit.Range("F47").Value = 448
$wsExhibit.Range("F49").Value = 183

$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F17").Value = 181
$wsShow.Range("F20").Value = 72
$wsShow.Range("F22").Value = 30

$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 43

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F10").Value = 9150
$wsAll.Range("F11").Value = 9150
$wsAll.Range("F12").Value = 535
$wsAll.Range("F15").Value = 220
$wsAll.Range("F17").Value = 136
$wsAll.Range("F19").Value = 11704
$wsAll.Range("F20").Value = 11704
$wsAll.Range("F24").Value = 43
$wsAll.Range("F30").Value = 168
$wsAll.Range("F31").Value = 145
$wsAll.Range("F35").Value = 2091
$wsAll.Range("F36").Value = 57
$wsAll.Range("F38").Value = 2132
$wsAll.Range("F39").Value = 961
$wsAll.Range("F40").Value = 181
$wsAll.Range("F42").Value = 3589
$wsAll.Range("F44").Value = 72
$wsAll.Range("F45").Value = 1295
$wsAll.Range("F47").Value = 388
$wsAll.Range("F48").Value = 30
$wsAll.Range("F49").Value = 448
$wsAll.Range("F51").Value = 183
